$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear A1 (it no longer holds "team1" - the header row has no label in col A)
$ws.Range("A1").Value = ""

# Row 1 becomes the header row: Nom / Score
$ws.Range("B1").Value = "Nom"
$ws.Range("C1").Value = "Score"

# Row 2: new "team info" row
$ws.Range("A2").Value = "Information de la première team"
$ws.Range("B2").Value = "VP"
$ws.Range("C2").Value = 5

# Row 3: previous team2 row, now shifted down
$ws.Range("A3").Value = "team2"
$ws.Range("B3").Value = "Astralis"
$ws.Range("C3").Value = 3

# Column A is widened to fit the long "team info" label
$ws.Columns.Item(1).ColumnWidth = 54.43

# Move the active selection to A5 (cursor resting below the data, unselected)
$null = $ws.Range("A5").Select()
